$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows of face/name data (person6 -> Malinga, person8 -> Hasith)
# Set column B (Image_Path) before column A (Name) on row 7 so that the
# shared-string table indices come out in the same order as the target file.
$ws.Range("B7").Value = "person6.jpg"
$ws.Range("A7").Value = "Malinga"
$ws.Range("A8").Value = "Hasith"
$ws.Range("B8").Value = "person8.jpg"

# Move the active selection to match the saved view state of the edited file
$ws.Range("J19").Select()

Write-Output "done"
